$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 16 and 17 text (wording tweaks)
$ws.Range("A16").Value = "Okay. Now could you please provide preferred date, required service and contact number?"
$ws.Range("A17").Value = "Okay Thank you CALLER! Now I'm going to make an appointment inquiry for USER for GSERVICE on APPDATE and contact number is PHONUM. Shall I confirm this?"

# New rows 18-22 completing the "Make Reservation" conversation flow
$ws.Range("A18").Value = "Thank you CALLER for the confirmation. Our staff member will check the availabilty and get back to as soon as possible. :) "
$ws.Range("B18").Value = 53

$ws.Range("A19").Value = "Sorry CALLER! Could you please tell me again the preferred  date, required service and contact number?"
$ws.Range("B19").Value = 54

$ws.Range("A20").Value = "Sorry, Phone number is not detected! Please type at least one contact number."
$ws.Range("B20").Value = 55

$ws.Range("A21").Value = "Sorry, Required service is not detected! Please type a service which is availble on our saloon."
$ws.Range("B21").Value = 56

$ws.Range("A22").Value = "Sorry, Preference date is not detected! Please type at least one preference date."
$ws.Range("B22").Value = 57

# Match style (wrap text) used by the other "Chat" column entries
$ws.Range("A18:A22").WrapText = $true

# Row heights (best-effort match of authored sizes)
$ws.Rows.Item(16).RowHeight = 32
$ws.Rows.Item(17).RowHeight = 51
$ws.Rows.Item(18).RowHeight = 29
$ws.Rows.Item(19).RowHeight = 34
$ws.Rows.Item(20).RowHeight = 17
$ws.Rows.Item(21).RowHeight = 34
$ws.Rows.Item(22).RowHeight = 17

# Selection / view state, matching the final saved state
$ws.Range("G22").Select()
